# Generate Report for Handoff
# Updates the localization-status workbook to reflect that "b.md" has now
# been handed off for zh-cn and de-de, with a new handback version warning.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3) moves from "Handed back: in sync
# with en-US" to "Ready for handoff" for both zh-cn (E) and de-de (F)
# columns, with an updated "Latest HO Xliff Generate Date" (G).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-12 04:39:14"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets its own handback details populated,
# mirroring what row 2 (a.md) already had, but with the new file/status.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-12 04:39:09"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/8845ca77ceb16bf33ad29e5ee4aeb40669cfcd91/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/fdea83c245e6cf1d0cf26be45ab380ffeee8b9c1/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: same kind of update as zh-cn, but with de-de file names
# and its own "Latest Handoff Datetime".
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-12 04:39:14"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/8845ca77ceb16bf33ad29e5ee4aeb40669cfcd91/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/fdea83c245e6cf1d0cf26be45ab380ffeee8b9c1/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
